# Weekly NYPD CompStat update: refresh report header (volume/week) and all
# crime-complaint figures for the new reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: issue number and covered date range ---
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# --- Crime complaint figures (Week to Date / 28 Day / Year to Date / 2 Year) ---
# Row 14: Murder
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 57.142857142857
$ws.Range("I14").Value = 98
$ws.Range("J14").Value = 102
$ws.Range("K14").Value = -3.92156862745
$ws.Range("L14").Value = -17.647058823529
$ws.Range("M14").Value = 1.030927835051
$ws.Range("N14").Value = -73.655913978494

# Row 15: Rape
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -33.333333333333
$ws.Range("F15").Value = 23
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = -17.857142857142
$ws.Range("I15").Value = 281
$ws.Range("J15").Value = 293
$ws.Range("K15").Value = -4.095563139931
$ws.Range("L15").Value = 6.844106463878
$ws.Range("M15").Value = 24.888888888888
$ws.Range("N15").Value = -46.679316888045

# Row 16: Robbery
$ws.Range("C16").Value = 84
$ws.Range("D16").Value = 88
$ws.Range("E16").Value = -4.545454545454
$ws.Range("F16").Value = 382
$ws.Range("G16").Value = 383
$ws.Range("H16").Value = -0.261096605744
$ws.Range("I16").Value = 3546
$ws.Range("J16").Value = 3741
$ws.Range("K16").Value = -5.212510024057
$ws.Range("L16").Value = 32.017870439315
$ws.Range("M16").Value = 11.194731890874
$ws.Range("N16").Value = -69.77239792004

# Row 17: Fel. Assault
$ws.Range("C17").Value = 137
$ws.Range("D17").Value = 150
$ws.Range("E17").Value = -8.666666666666
$ws.Range("F17").Value = 677
$ws.Range("G17").Value = 611
$ws.Range("H17").Value = 10.801963993453
$ws.Range("I17").Value = 5991
$ws.Range("J17").Value = 5446
$ws.Range("K17").Value = 10.007344840249
$ws.Range("L17").Value = 32.661647475642
$ws.Range("M17").Value = 81.106408706166
$ws.Range("N17").Value = -11.545843791525

# Row 18: Burglary
$ws.Range("C18").Value = 54
$ws.Range("D18").Value = 49
$ws.Range("E18").Value = 10.204081632653
$ws.Range("F18").Value = 219
$ws.Range("G18").Value = 244
$ws.Range("H18").Value = -10.245901639344
$ws.Range("I18").Value = 2185
$ws.Range("J18").Value = 2154
$ws.Range("K18").Value = 1.439182915506
$ws.Range("L18").Value = 39.171974522293
$ws.Range("M18").Value = -8.57740585774
$ws.Range("N18").Value = -84.164371648064

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 150
$ws.Range("D19").Value = 149
$ws.Range("E19").Value = 0.671140939597
$ws.Range("F19").Value = 658
$ws.Range("G19").Value = 629
$ws.Range("H19").Value = 4.610492845786
$ws.Range("I19").Value = 5831
$ws.Range("J19").Value = 5840
$ws.Range("K19").Value = -0.154109589041
$ws.Range("L19").Value = 24.620645437059
$ws.Range("M19").Value = 70.397428404441
$ws.Range("N19").Value = 4.704614832106

# Row 20: G.L.A.
$ws.Range("C20").Value = 97
$ws.Range("D20").Value = 77
$ws.Range("E20").Value = 25.974025974026
$ws.Range("F20").Value = 408
$ws.Range("G20").Value = 278
$ws.Range("H20").Value = 46.762589928057
$ws.Range("I20").Value = 3925
$ws.Range("J20").Value = 2847
$ws.Range("K20").Value = 37.864418686336
$ws.Range("L20").Value = 91.744015632633
$ws.Range("M20").Value = 157.039947609692
$ws.Range("N20").Value = -64.989742217465

# Row 21: TOTAL
$ws.Range("C21").Value = 527
$ws.Range("D21").Value = 521
$ws.Range("E21").Value = 1.151631477927
$ws.Range("F21").Value = 2378
$ws.Range("G21").Value = 2180
$ws.Range("H21").Value = 9.082568807339
$ws.Range("I21").Value = 21857
$ws.Range("J21").Value = 20423
$ws.Range("K21").Value = 7.021495372863
$ws.Range("L21").Value = 37.638539042821
$ws.Range("M21").Value = 54.3791495974
$ws.Range("N21").Value = -56.269382365298

# Row 22: Transit
$ws.Range("C22").Value = 5
$ws.Range("E22").Value = -37.5
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 34
$ws.Range("H22").Value = -23.529411764705
$ws.Range("I22").Value = 213
$ws.Range("J22").Value = 269
$ws.Range("K22").Value = -20.817843866171
$ws.Range("L22").Value = 18.994413407821
$ws.Range("M22").Value = -6.986899563318

# Row 23: Housing
$ws.Range("C23").Value = 23
$ws.Range("D23").Value = 33
$ws.Range("E23").Value = -30.30303030303
$ws.Range("F23").Value = 123
$ws.Range("G23").Value = 133
$ws.Range("H23").Value = -7.518796992481
$ws.Range("I23").Value = 1308
$ws.Range("J23").Value = 1183
$ws.Range("K23").Value = 10.566356720202
$ws.Range("L23").Value = 49.31506849315
$ws.Range("M23").Value = 66.624203821656

# Row 24: Petit Larceny
$ws.Range("C24").Value = 378
$ws.Range("D24").Value = 380
$ws.Range("E24").Value = -0.526315789473
$ws.Range("F24").Value = 1449
$ws.Range("G24").Value = 1465
$ws.Range("H24").Value = -1.092150170648
$ws.Range("I24").Value = 13203
$ws.Range("J24").Value = 13609
$ws.Range("K24").Value = -2.983319861856
$ws.Range("L24").Value = 42.304375943091
$ws.Range("M24").Value = 41.042623651319

# Row 25: Misd. Assault
$ws.Range("C25").Value = 178
$ws.Range("D25").Value = 222
$ws.Range("E25").Value = -19.819819819819
$ws.Range("F25").Value = 783
$ws.Range("G25").Value = 765
$ws.Range("H25").Value = 2.35294117647
$ws.Range("I25").Value = 7715
$ws.Range("J25").Value = 7388
$ws.Range("K25").Value = 4.426096372495
$ws.Range("L25").Value = 23.697290363957
$ws.Range("M25").Value = -6.314511232544

# Row 26: UCR Rape*
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -26.666666666666
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = -12
$ws.Range("I26").Value = 480
$ws.Range("J26").Value = 510
$ws.Range("K26").Value = -5.882352941176
$ws.Range("L26").Value = 11.111111111111

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 29
$ws.Range("D27").Value = 23
$ws.Range("E27").Value = 26.086956521739
$ws.Range("F27").Value = 102
$ws.Range("G27").Value = 73
$ws.Range("H27").Value = 39.72602739726
$ws.Range("I27").Value = 775
$ws.Range("J27").Value = 673
$ws.Range("K27").Value = 15.156017830609
$ws.Range("L27").Value = 16.541353383458

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = 14.285714285714
$ws.Range("G28").Value = 37
$ws.Range("H28").Value = -10.81081081081
$ws.Range("I28").Value = 305
$ws.Range("J28").Value = 384
$ws.Range("K28").Value = -20.572916666666
$ws.Range("L28").Value = -32.819383259911
$ws.Range("M28").Value = -17.119565217391
$ws.Range("N28").Value = -70.896946564885

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = -14.285714285714
$ws.Range("F29").Value = 27
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = -3.571428571428
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 323
$ws.Range("K29").Value = -22.600619195046
$ws.Range("L29").Value = -34.383202099737
$ws.Range("M29").Value = -18.566775244299
$ws.Range("N29").Value = -73.544973544973

# Row 30: Hate Crimes
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 18
$ws.Range("K30").Value = -47.058823529411
$ws.Range("L30").Value = -48.571428571428
